# routes.xlsx update:
#   - append a new route row (East King County / 98288 / east_king_county)
#   - grow the "Routes" query table + used range to include it
#   - keep the ExternalData_1 defined name (used by the Power Query
#     refresh) pointing at the full, newly-sized range
#   - move the visible selection down near the new row, like a user
#     would leave it after typing the new data in at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("routes")
$ws.Activate()

$tbl = $ws.ListObjects.Item("Routes")

# Grow the table by one row, then fill in the new record.
$newListRow = $tbl.ListRows.Add()
$newRow = $newListRow.Range.Row

$ws.Cells.Item($newRow, 1).Value = "East King County"
$ws.Cells.Item($newRow, 2).Value = 98288
$ws.Cells.Item($newRow, 3).Value = "east_king_county"

# The workbook's ExternalData_1 name tracks the query table's full range
# (it is what Excel uses to remember the last-refreshed extent) - keep it
# in sync with the resized table.
$extName = $wb.Names.Item("ExternalData_1")
$extName.RefersTo = "=routes!" + $tbl.Range.Address()

# Leave the view scrolled near the bottom, with the last-edited cell
# selected, matching where a person would be after adding this row.
$ws.Range("C96").Select()
